$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Costos")

# Update the "Licencias de programas" cost from 5000 to 25000
$ws.Range("C16").Value = 25000

# Update the client's budget from 500000 to 450000
$ws.Range("D40").Value = 450000

# Update the selection/view state on the Costos sheet
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("D34").Select()
